# Apply updated crypto price/volume data (Thu Nov 28 19:31:19 UTC 2024 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force text storage (matches source inlineStr cells) without leaving a
    # residual number-format style stamp on the cell once the write lands.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "95.092.99"
Set-TextValue $ws.Range("E2") "  -1.30%  "
Set-TextValue $ws.Range("D3") "3.569.29"
Set-TextValue $ws.Range("E3") "  -1.15%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "235.45"
Set-TextValue $ws.Range("E5") "  -1.99%  "
Set-TextValue $ws.Range("D6") "653.28"
Set-TextValue $ws.Range("E6") "  +2.33%  "
Set-TextValue $ws.Range("E7") "  -0.74%  "
Set-TextValue $ws.Range("E8") "  -0.93%  "
Set-TextValue $ws.Range("E9") "  +0.14%  "
Set-TextValue $ws.Range("E10") "  -2.04%  "
Set-TextValue $ws.Range("D11") "3.566.05"
Set-TextValue $ws.Range("E11") "  -1.19%  "
Set-TextValue $ws.Range("E12") "  +1.10%  "
Set-TextValue $ws.Range("D13") "42.39"
Set-TextValue $ws.Range("E13") "  -1.89%  "
Set-TextValue $ws.Range("D14") "6.49"
Set-TextValue $ws.Range("E14") "  +1.27%  "
Set-TextValue $ws.Range("D15") "4.273.40"
Set-TextValue $ws.Range("E15") "  -0.38%  "
Set-TextValue $ws.Range("D16") "95.029.26"
Set-TextValue $ws.Range("E16") "  -1.30%  "
Set-TextValue $ws.Range("E17") "  -0.31%  "
Set-TextValue $ws.Range("D18") "3.569.69"
Set-TextValue $ws.Range("E18") "  -1.00%  "
Set-TextValue $ws.Range("D19") "7.74"
Set-TextValue $ws.Range("E19") "  -7.55%  "
Set-TextValue $ws.Range("D20") "12.58"
Set-TextValue $ws.Range("E20") "  -5.06%  "
Set-TextValue $ws.Range("D21") "17.83"
Set-TextValue $ws.Range("E21") "  -1.65%  "
Set-TextValue $ws.Range("E22") "  -0.17%  "
Set-TextValue $ws.Range("D23") "508.48"
Set-TextValue $ws.Range("E23") "  -1.50%  "
Set-TextValue $ws.Range("E24") "  -4.50%  "
Set-TextValue $ws.Range("D25") "6.78"
Set-TextValue $ws.Range("E25") "  +1.76%  "
Set-TextValue $ws.Range("E26") "  -1.91%  "
Set-TextValue $ws.Range("D27") "95.24"
Set-TextValue $ws.Range("E27") "  -1.88%  "
Set-TextValue $ws.Range("D28") "12.73"
Set-TextValue $ws.Range("E28") "  +2.17%  "
Set-TextValue $ws.Range("D29") "3.760.85"
Set-TextValue $ws.Range("E29") "  -0.99%  "
Set-TextValue $ws.Range("D30") "3.04"
Set-TextValue $ws.Range("E30") "  -1.40%  "
Set-TextValue $ws.Range("D31") "11.52"
Set-TextValue $ws.Range("E31") "  -1.43%  "
Set-TextValue $ws.Range("E32") "  -0.10%  "
Set-TextValue $ws.Range("E33") "  +0.06%  "
Set-TextValue $ws.Range("E34") "  +0.73%  "
Set-TextValue $ws.Range("E35") "  -2.15%  "
Set-TextValue $ws.Range("D36") "31.71"
Set-TextValue $ws.Range("E36") "  +4.00%  "
Set-TextValue $ws.Range("D37") "1.69"
Set-TextValue $ws.Range("E37") "  +13.00%  "
Set-TextValue $ws.Range("D38") "8.54"
Set-TextValue $ws.Range("E38") "  +8.33%  "
Set-TextValue $ws.Range("E39") "  -1.81%  "
Set-TextValue $ws.Range("D40") "583.38"
Set-TextValue $ws.Range("E40") "  +1.33%  "
Set-TextValue $ws.Range("E41") "  +0.07%  "
Set-TextValue $ws.Range("E42") "  -0.92%  "
Set-TextValue $ws.Range("D43") "0.904"
Set-TextValue $ws.Range("E43") "  -2.82%  "
Set-TextValue $ws.Range("D44") "1.82"
Set-TextValue $ws.Range("E44") "  +3.66%  "
Set-TextValue $ws.Range("E45") "  +4.55%  "
Set-TextValue $ws.Range("E46") "  +1.14%  "
Set-TextValue $ws.Range("D47") "34.04"
Set-TextValue $ws.Range("E47") "  +29.62%  "
Set-TextValue $ws.Range("E48") "  -1.75%  "
Set-TextValue $ws.Range("E49") "  -3.59%  "
Set-TextValue $ws.Range("E50") "  +0.67%  "
Set-TextValue $ws.Range("D51") "8.18"
Set-TextValue $ws.Range("E51") "  +0.19%  "
